$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format columns D and F (fit_time / predict_time) as Text so the
# numeric-looking strings ("1.132", "0.045", ...) are stored verbatim
# instead of being coerced to floating point numbers by Excel.
$ws.Range("D2:D21").NumberFormat = "@"
$ws.Range("F2:F21").NumberFormat = "@"

$data = @(
    @("christine", "xgb_baseline", "(3467, 1636)", "1.132", "(1084, 1636)", "0.045", 0.7306273062730627, 0.7306638344515817, 0.7310098423921751, 0.7255639097744362, 0.7933523794588219),
    @("christine", "lgbm_baseline", "(3467, 1636)", "8.614", "(1084, 1636)", "0.026", 0.7416974169741697, 0.7413716173640559, 0.7415954425667268, 0.7333333333333332, 0.8270677717055941),
    @("christine", "rf_baseline", "(3467, 1636)", "1.871", "(1084, 1636)", "0.023", 0.7214022140221402, 0.7213017650763391, 0.7216276937028463, 0.7150943396226415, 0.7938620565042632),
    @("christine", "logit_baseline", "(3467, 1636)", "132.973", "(1084, 1636)", "0.009", 0.6171586715867159, 0.6170220733417759, 0.6171916283414529, 0.6081208687440982, 0.6671382741892045),
    @("jasmine", "xgb_baseline", "(1909, 144)", "0.038", "(597, 144)", "0.004", 0.8107202680067002, 0.830954954954955, 0.8092602985744752, 0.8330871491875923, 0.8675103827590077),
    @("jasmine", "lgbm_baseline", "(1909, 144)", "0.233", "(597, 144)", "0.007", 0.8040201005025126, 0.814968487394958, 0.8029127848243349, 0.8224582701062215, 0.8710854192389719),
    @("jasmine", "rf_baseline", "(1909, 144)", "0.253", "(597, 144)", "0.015", 0.7839195979899497, 0.8407516100105976, 0.7815130766640477, 0.8220689655172414, 0.8644404534740149),
    @("jasmine", "logit_baseline", "(1909, 144)", "0.024", "(597, 144)", "0.003", 0.7772194304857621, 0.7850045167118338, 0.7762262880233471, 0.7963246554364471, 0.8537770793579526),
    @("philippine", "xgb_baseline", "(3732, 308)", "0.346", "(1167, 308)", "0.013", 0.7557840616966581, 0.7596903129462936, 0.7577429334337101, 0.7484554280670784, 0.8453104231117593),
    @("philippine", "lgbm_baseline", "(3732, 308)", "1.735", "(1167, 308)", "0.010", 0.7660668380462725, 0.7682151913663974, 0.7674958535753356, 0.7628149435273676, 0.8611212402808982),
    @("philippine", "rf_baseline", "(3732, 308)", "1.393", "(1167, 308)", "0.020", 0.7052270779777207, 0.7320433316287229, 0.7108309905543858, 0.6594059405940594, 0.801292155317422),
    @("philippine", "logit_baseline", "(3732, 308)", "3.765", "(1167, 308)", "0.004", 0.7300771208226221, 0.7339717741935483, 0.7320645077811629, 0.7214854111405835, 0.7972751226283039),
    @("madeline", "xgb_baseline", "(2009, 259)", "0.097", "(628, 259)", "0.007", 0.8264331210191083, 0.8263654343526547, 0.8270121951219511, 0.8304821150855364, 0.8932113821138212),
    @("madeline", "lgbm_baseline", "(2009, 259)", "0.781", "(628, 259)", "0.010", 0.8152866242038217, 0.8152866242038217, 0.8159146341463415, 0.8193146417445484, 0.8934247967479675),
    @("madeline", "rf_baseline", "(2009, 259)", "0.532", "(628, 259)", "0.017", 0.6385350318471338, 0.6391011395577744, 0.6393089430894309, 0.6425196850393701, 0.7398069105691056),
    @("madeline", "logit_baseline", "(2009, 259)", "0.653", "(628, 259)", "0.003", 0.5700636942675159, 0.5694241390513319, 0.5694918699186992, 0.5858895705521472, 0.598709349593496),
    @("sylvine", "xgb_baseline", "(3279, 20)", "0.034", "(1025, 20)", "0.003", 0.9395121951219512, 0.9422083735022724, 0.9402077049629832, 0.94106463878327, 0.9770397282393444),
    @("sylvine", "lgbm_baseline", "(3279, 20)", "0.184", "(1025, 20)", "0.006", 0.9414634146341463, 0.9436491072520484, 0.9420947201657375, 0.9427480916030535, 0.9827026780001828),
    @("sylvine", "rf_baseline", "(3279, 20)", "0.397", "(1025, 20)", "0.017", 0.9063414634146342, 0.9103117588819046, 0.9071896992962252, 0.9096045197740114, 0.9556332449806538),
    @("sylvine", "logit_baseline", "(3279, 20)", "0.085", "(1025, 20)", "0.004", 0.9004878048780488, 0.9029958377327886, 0.9011725771562624, 0.903041825095057, 0.9598223806477165),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = [double]$row[6]
    $ws.Cells.Item($r, 8).Value = [double]$row[7]
    $ws.Cells.Item($r, 9).Value = [double]$row[8]
    $ws.Cells.Item($r, 10).Value = [double]$row[9]
    $ws.Cells.Item($r, 11).Value = [double]$row[10]
}

Write-Output "done"